# Applies the diff to Saldo.xlsx / Export sheet:
#  - Delete the old "001922009 / SOFIA / 81.11" row (Excel row 180)
#  - Delete the "004504449 / KELMA / 1000" row (Excel row 5)
#  - Insert a new row before "004466342 / TATYANA" (Excel row 4) with
#    "001922009 / SOFIA / 10514.44"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old SOFIA row (row 180), bottom-most change first so the
#    row numbers used below remain valid.
$ws.Rows.Item(180).Delete()

# 2) Delete the KELMA row (row 5).
$ws.Rows.Item(5).Delete()

# 3) Insert a new row before row 4 (TATYANA) and populate it with the
#    relocated SOFIA record (now with a different balance).
$ws.Rows.Item(4).Insert()

# Preserve the leading zeros of the account number by forcing the cell
# to be text before assigning the value.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "001922009"
$ws.Range("B4").Value = "SOFIA"
$ws.Range("C4").Value = 10514.44

Write-Host ("A4=" + $ws.Range("A4").Value() + " B4=" + $ws.Range("B4").Value() + " C4=" + $ws.Range("C4").Value())
Write-Host ("Row5 A=" + $ws.Range("A5").Value() + " B5=" + $ws.Range("B5").Value())
Write-Host ("UsedRangeRows=" + $ws.UsedRange.Rows.Count)
